$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Startseite: re-order the "Eigenfertigung" block and add the new
#    "Terrassenüberdachung" entry.
#
#    Old rows 17-18:
#      17: >> Montagematerial | Zub_Montage
#      18: Glasdach            | Eigen_Glasdach
#
#    New rows 17-19:
#      17: Glasdach              | Eigen_Glasdach
#      18: Terrassenüberdachung  | Eigen_Terrasse
#      19: >> Montagematerial    | Zub_Montage
# ---------------------------------------------------------------------------
$start = $wb.Worksheets.Item(1)

$start.Cells.Item(17, 2).Value = "Glasdach"
$start.Cells.Item(17, 3).Value = "Eigen_Glasdach"

$start.Cells.Item(18, 2).Value = "Terrassenüberdachung"
$start.Cells.Item(18, 3).Value = "Eigen_Terrasse"

$start.Cells.Item(19, 1).Value = "Eigenfertigung"
$start.Cells.Item(19, 2).Value = ">> Montagematerial"
$start.Cells.Item(19, 3).Value = "Zub_Montage"

$start.Columns.Item(2).ColumnWidth = 40.83333333333333
$start.Range("C18").Select()

# ---------------------------------------------------------------------------
# 2) Insert the new "Eigen_Terrasse" worksheet right after "Startseite".
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "Eigen_Terrasse"

$ws.Cells.Item(1, 1).Value = "Zeile"
$ws.Cells.Item(1, 2).Value = "Typ"
$ws.Cells.Item(1, 3).Value = "Bezeichnung"
$ws.Cells.Item(1, 4).Value = "Variable"
$ws.Cells.Item(1, 5).Value = "Optionen (Beispielwerte)"
$ws.Cells.Item(1, 6).Value = "Formel (Logik)"

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Zahl"
$ws.Cells.Item(2, 3).Value = "Länge (m)"
$ws.Cells.Item(2, 4).Value = "L"

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Zahl"
$ws.Cells.Item(3, 3).Value = "Tiefe/Breite (m)"
$ws.Cells.Item(3, 4).Value = "B"

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Zahl"
$ws.Cells.Item(4, 3).Value = "Säulen Höhe (m)"
$ws.Cells.Item(4, 4).Value = "H"
$ws.Cells.Item(4, 5).Value = 2.5
$ws.Cells.Item(4, 5).NumberFormat = "0.00"

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Zahl"
$ws.Cells.Item(5, 3).Value = "Anzahl Säulen"
$ws.Cells.Item(5, 4).Value = "N_Col"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 5).NumberFormat = "0.00"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Auswahl"
$ws.Cells.Item(6, 3).Value = "Dach-Eindeckung (€/m²)"
$ws.Cells.Item(6, 4).Value = "P_Dach"
$ws.Cells.Item(6, 5).Value = "VSG Glas 10mm:180, Flachdach Folie:140"
$ws.Cells.Item(6, 6).Value = "Preis inkl. Gummi/Leisten"

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Zahl"
$ws.Cells.Item(7, 3).Value = "Wandanschluss (€/m)"
$ws.Cells.Item(7, 4).Value = "P_Wand"
$ws.Cells.Item(7, 5).Value = 75
$ws.Cells.Item(7, 6).Value = "Blech & Flüssigkunststoff"

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Zahl"
$ws.Cells.Item(8, 3).Value = "Rabatt (%)"
$ws.Cells.Item(8, 4).Value = "Rabatt"
$ws.Cells.Item(8, 5).Value = 0

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Berechnung"
$ws.Cells.Item(9, 3).Value = "Intern: Sparren-Anzahl"
$ws.Cells.Item(9, 4).Value = "N_Spar"
$ws.Cells.Item(9, 6).Value = "math.ceil(L / 1.0) + 1"

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Berechnung"
$ws.Cells.Item(10, 3).Value = "Intern: Trägerpreis/m"
$ws.Cells.Item(10, 4).Value = "P_Trager"
$ws.Cells.Item(10, 6).Value = "(L <= 5) * 60 + (L > 5) * 110"

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Preis"
$ws.Cells.Item(11, 3).Value = "Gesamtpreis"
$ws.Cells.Item(11, 4).Value = "Endpreis"
$ws.Cells.Item(11, 6).Value = "( (L * P_Trager) + (N_Col * H * 50) + (N_Spar * B * 45) + (L * B * P_Dach) + (L * P_Wand) ) * (1 - (Rabatt / 100))"

$ws.Columns.Item(1).ColumnWidth = 4.666666666666666
$ws.Columns.Item(2).ColumnWidth = 10.666666666666668
$ws.Columns.Item(3).ColumnWidth = 21.666666666666664
$ws.Columns.Item(4).ColumnWidth = 7.833333333333334
$ws.Columns.Item(5).ColumnWidth = 36.16666666666667
$ws.Columns.Item(6).ColumnWidth = 23.333333333333336

$ws.Activate()
$ws.Range("F20").Select()
